$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column C for rows 2-6
# from serial date 45224 (2023-10-25) to 45233 (2023-11-03)
foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = 45233
}
